$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a plain (non date-like) value into a cell so that it inherits
# the worksheet column's default style (matches how Excel behaves when a
# brand-new value is typed into a cell that has no direct cell-level style
# override yet).
# ---------------------------------------------------------------------------
function Set-PlainValue($ws, $addr, $value) {
    $ws.Range($addr).ClearContents()
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------------
# Helper: write literal text that Excel would otherwise auto-convert (e.g.
# dates like "2021-10-06") while still keeping the column's default style
# and the plain "text" cell type. We stage the text in a distant scratch
# cell formatted as Text (so Excel does not reinterpret it), copy only the
# resulting value (not the scratch formatting) onto the destination, then
# clean the scratch cell back up.
# ---------------------------------------------------------------------------
function Set-LiteralText($ws, $addr, $text) {
    $ws.Range($addr).ClearContents()
    $ws.Range($addr).Value = "x"

    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Same idea as Set-LiteralText, but for a destination cell that must keep no
# explicit style at all (i.e. it should render exactly like $styleSource).
# The cell's number format/font/etc. is first copied (format-only) from
# $styleSource before the literal text value is dropped in.
# ---------------------------------------------------------------------------
function Set-LiteralTextWithStyleFrom($ws, $addr, $text, $styleSource) {
    $styleSource.Copy()
    $ws.Range($addr).PasteSpecial(-4122)

    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Helper: write a date/time serial number into a cell, forcing the
# "YYYY-MM-DD HH:MM:SS" custom number format (re-using the workbook's
# existing number format definition) regardless of the cell's prior format.
# ---------------------------------------------------------------------------
function Set-DateTimeValue($ws, $addr, $serial) {
    $ws.Range($addr).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range($addr).Value = $serial
}

# ---------------------------------------------------------------------------
# Helper: plain value, but with the cell first force-styled (format only)
# from a reference cell so it ends up with no direct style override.
# ---------------------------------------------------------------------------
function Set-PlainValueWithStyleFrom($ws, $addr, $value, $styleSource) {
    $styleSource.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $value
}

$ws1 = $wb.Worksheets.Item("AMSIN")
$ws2 = $wb.Worksheets.Item("BETA")
$ws3 = $wb.Worksheets.Item("AMS")

# A cell that (at this point, before any edits) has no explicit style of its
# own - used as a formatting donor for the "AMS" sheet's new row further
# below, which must come out with no direct style override either.
$noStyleDonor = $ws1.Range("D23")

# ===========================================================================
# Sheet "AMS": row 18 is new - its cells must end up without any explicit
# style override (same rendering as the workbook default), so this is done
# first, while $noStyleDonor is still unstyled.
# ===========================================================================
Set-LiteralTextWithStyleFrom $ws3 "A18" "2021-10-28" $noStyleDonor
Set-DateTimeValue $ws3 "B18" 44497.85925816281
Set-PlainValueWithStyleFrom $ws3 "C18" "152_livetest" $noStyleDonor
Set-PlainValue $ws3 "D18" 124
Set-PlainValueWithStyleFrom $ws3 "E18" 118 $noStyleDonor
Set-PlainValueWithStyleFrom $ws3 "F18" 6 $noStyleDonor
Set-PlainValueWithStyleFrom $ws3 "G18" 2.56 $noStyleDonor

# ===========================================================================
# Sheet "AMSIN": row 23 gets restyled (same values), rows 24-26 are new.
# ===========================================================================

# Row 23 - values stay the same, only the cell styling is normalized to
# match the rest of the table, and the run-time float is refreshed.
Set-LiteralText $ws1 "A23" "2021-10-06"
Set-DateTimeValue $ws1 "B23" 44475.70126787037
Set-PlainValue $ws1 "C23" "151_regression"
Set-PlainValue $ws1 "D23" 124
Set-PlainValue $ws1 "E23" 120
Set-PlainValue $ws1 "F23" 4
Set-PlainValue $ws1 "G23" 2.35

# Row 24 - new sprint entry.
Set-LiteralText $ws1 "A24" "2021-10-26"
Set-DateTimeValue $ws1 "B24" 44495.64759021991
Set-PlainValue $ws1 "C24" "152_fstcycle"
Set-PlainValue $ws1 "D24" 124
Set-PlainValue $ws1 "E24" 122
Set-PlainValue $ws1 "F24" 2
Set-PlainValue $ws1 "G24" 2.37

# Row 25 - new sprint entry.
Set-LiteralText $ws1 "A25" "2021-10-27"
Set-DateTimeValue $ws1 "B25" 44496.64631342592
Set-PlainValue $ws1 "C25" "152_scndcycle"
Set-PlainValue $ws1 "D25" 124
Set-PlainValue $ws1 "E25" 120
Set-PlainValue $ws1 "F25" 4
Set-PlainValue $ws1 "G25" 3.55

# Row 26 - new sprint entry.
Set-LiteralText $ws1 "A26" "2021-10-28"
Set-DateTimeValue $ws1 "B26" 44497.38627040509
Set-PlainValue $ws1 "C26" "152_fnlrgrsn"
Set-PlainValue $ws1 "D26" 124
Set-PlainValue $ws1 "E26" 115
Set-PlainValue $ws1 "F26" 9
Set-PlainValue $ws1 "G26" 3.02

# ===========================================================================
# Sheet "BETA": row 15 is new.
# ===========================================================================

Set-LiteralText $ws2 "A15" "2021-10-28"
Set-DateTimeValue $ws2 "B15" 44497.62344630787
Set-PlainValue $ws2 "C15" "152_beta"
Set-PlainValue $ws2 "D15" 124
Set-PlainValue $ws2 "E15" 117
Set-PlainValue $ws2 "F15" 7
Set-PlainValue $ws2 "G15" 3.03
